$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.063713
$ws.Range("H2").Value = 12.191139
$ws.Range("I2").Value = 0.5065008440615062
$ws.Range("J2").Value = 0.5065008440615063
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07289133333333334
$ws.Range("N2").Value = 0.218674
$ws.Range("O2").Value = 0.0516105978808581
$ws.Range("P2").Value = 0.0516105978808581
$ws.Range("Q2").Value = 0.296209458854
$ws.Range("R2").Value = 2.665885129686
$ws.Range("S2").Value = 0.02614081138917361
$ws.Range("T2").Value = 0.02614081138917362
$ws.Range("G3").Value = 4.063713
$ws.Range("H3").Value = 12.191139
$ws.Range("I3").Value = 0.5065008440615062
$ws.Range("J3").Value = 0.5065008440615063
$ws.Range("O3").Value = 0.4649150176610893
$ws.Range("P3").Value = 0.4649150176610893
$ws.Range("Q3").Value = 2.668293556924
$ws.Range("R3").Value = 24.014642012316
$ws.Range("S3").Value = 0.2354798488622118
$ws.Range("T3").Value = 0.2354798488622118
$ws.Range("G4").Value = 4.063713
$ws.Range("H4").Value = 12.191139
$ws.Range("I4").Value = 0.5065008440615062
$ws.Range("J4").Value = 0.5065008440615063
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6828266666666667
$ws.Range("N4").Value = 2.04848
$ws.Range("O4").Value = 0.4834743844580526
$ws.Range("P4").Value = 0.4834743844580526
$ws.Range("Q4").Value = 2.77481160208
$ws.Range("R4").Value = 24.97330441872
$ws.Range("S4").Value = 0.2448801838101208
$ws.Range("T4").Value = 0.2448801838101209
$ws.Range("I5").Value = 0.2604012840237886
$ws.Range("J5").Value = 0.2604012840237886
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07289133333333334
$ws.Range("N5").Value = 0.218674
$ws.Range("O5").Value = 0.0516105978808581
$ws.Range("P5").Value = 0.0516105978808581
$ws.Range("Q5").Value = 0.1522866631515556
$ws.Range("R5").Value = 1.370579968364
$ws.Range("S5").Value = 0.01343946595741087
$ws.Range("T5").Value = 0.01343946595741088
$ws.Range("I6").Value = 0.2604012840237886
$ws.Range("J6").Value = 0.2604012840237886
$ws.Range("O6").Value = 0.4649150176610893
$ws.Range("P6").Value = 0.4649150176610893
$ws.Range("S6").Value = 0.12106446756089
$ws.Range("T6").Value = 0.12106446756089
$ws.Range("I7").Value = 0.2604012840237886
$ws.Range("J7").Value = 0.2604012840237886
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6828266666666667
$ws.Range("N7").Value = 2.04848
$ws.Range("O7").Value = 0.4834743844580526
$ws.Range("P7").Value = 0.4834743844580526
$ws.Range("Q7").Value = 1.426581046364444
$ws.Range("R7").Value = 12.83922941728
$ws.Range("S7").Value = 0.1258973505054877
$ws.Range("T7").Value = 0.1258973505054878
$ws.Range("G8").Value = 1.588356333333333
$ws.Range("H8").Value = 4.765069
$ws.Range("I8").Value = 0.1979725988286506
$ws.Range("J8").Value = 0.1979725988286507
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07289133333333334
$ws.Range("N8").Value = 0.218674
$ws.Range("O8").Value = 0.0516105978808581
$ws.Range("P8").Value = 0.0516105978808581
$ws.Range("Q8").Value = 0.1157774109451111
$ws.Range("R8").Value = 1.041996698506
$ws.Range("S8").Value = 0.01021748418957393
$ws.Range("T8").Value = 0.01021748418957393
$ws.Range("G9").Value = 1.588356333333333
$ws.Range("H9").Value = 4.765069
$ws.Range("I9").Value = 0.1979725988286506
$ws.Range("J9").Value = 0.1979725988286507
$ws.Range("O9").Value = 0.4649150176610893
$ws.Range("P9").Value = 0.4649150176610893
$ws.Range("Q9").Value = 1.042938064359556
$ws.Range("R9").Value = 9.386442579236
$ws.Range("S9").Value = 0.09204043428083386
$ws.Range("T9").Value = 0.09204043428083389
$ws.Range("G10").Value = 1.588356333333333
$ws.Range("H10").Value = 4.765069
$ws.Range("I10").Value = 0.1979725988286506
$ws.Range("J10").Value = 0.1979725988286507
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6828266666666667
$ws.Range("N10").Value = 2.04848
$ws.Range("O10").Value = 0.4834743844580526
$ws.Range("P10").Value = 0.4834743844580526
$ws.Range("Q10").Value = 1.084572060568889
$ws.Range("R10").Value = 9.761148545119999
$ws.Range("S10").Value = 0.09571468035824286
$ws.Range("T10").Value = 0.09571468035824288
$ws.Range("G11").Value = 0.281814
$ws.Range("H11").Value = 0.845442
$ws.Range("I11").Value = 0.03512527308605438
$ws.Range("J11").Value = 0.03512527308605439
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.07289133333333334
$ws.Range("N11").Value = 0.218674
$ws.Range("O11").Value = 0.0516105978808581
$ws.Range("P11").Value = 0.0516105978808581
$ws.Range("Q11").Value = 0.020541798212
$ws.Range("R11").Value = 0.184876183908
$ws.Range("S11").Value = 0.00181283634469968
$ws.Range("T11").Value = 0.001812836344699681
$ws.Range("G12").Value = 0.281814
$ws.Range("H12").Value = 0.845442
$ws.Range("I12").Value = 0.03512527308605438
$ws.Range("J12").Value = 0.03512527308605439
$ws.Range("O12").Value = 0.4649150176610893
$ws.Range("P12").Value = 0.4649150176610893
$ws.Range("Q12").Value = 0.185043205672
$ws.Range("R12").Value = 1.665388851048
$ws.Range("S12").Value = 0.01633026695715356
$ws.Range("T12").Value = 0.01633026695715356
$ws.Range("G13").Value = 0.281814
$ws.Range("H13").Value = 0.845442
$ws.Range("I13").Value = 0.03512527308605438
$ws.Range("J13").Value = 0.03512527308605439
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6828266666666667
$ws.Range("N13").Value = 2.04848
$ws.Range("O13").Value = 0.4834743844580526
$ws.Range("P13").Value = 0.4834743844580526
$ws.Range("Q13").Value = 0.19243011424
$ws.Range("R13").Value = 1.73187102816
$ws.Range("S13").Value = 0.01698216978420115
$ws.Range("T13").Value = 0.01698216978420115
